# Update "想去人数" (attendance interest count) figures in column F
# for rows 3-8 on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$newValues = @{
    3 = 2286
    4 = 386
    5 = 81
    6 = 6446
    7 = 316
    8 = 122
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $newValues[$row]
    }
}
